# Commit: re-run RU 1001; without crop
# The Russia (column L) responses for question 1001 were re-processed
# ("without crop"), which also changes the "All" aggregate column (B).
# A handful of answer-category rows had their label/value association
# corrected to line up with the recomputed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row label <-> data alignment ---
$ws.Range("A8").Value2 = "Security; violence; crime; judicial system"
$ws.Range("A9").Value2 = "Corruption; criticism of the government"
$ws.Range("A11").Value2 = "Nothing; don't know; empty"
$ws.Range("A12").Value2 = "Money; own income; cost of living; inflation"
$ws.Range("A17").Value2 = "Happiness; peace of mind"
$ws.Range("A18").Value2 = "Family; children; childcare"
$ws.Range("A19").Value2 = "Criticism of far right; Trump; tariffs"
$ws.Range("A20").Value2 = "Criticism of immigration; national preference"
$ws.Range("A21").Value2 = "Housing"
$ws.Range("A23").Value2 = "Social division; fake news; (social) media"
$ws.Range("A24").Value2 = "Education"

# --- Updated numeric values (re-run of RU / All aggregate, and rows whose data moved with their label) ---
$ws.Range("B2").Value2 = 0.273124301673209
$ws.Range("L2").Value2 = 0.160657026907098
$ws.Range("B3").Value2 = 0.21883268118649
$ws.Range("L3").Value2 = 0.204031066972921
$ws.Range("B4").Value2 = 0.204739119911163
$ws.Range("L4").Value2 = 0.0730309302608108
$ws.Range("B5").Value2 = 0.194193936161959
$ws.Range("L5").Value2 = 0.0231344891074306
$ws.Range("B6").Value2 = 0.160020718419108
$ws.Range("L6").Value2 = 0.0224941216710314
$ws.Range("B7").Value2 = 0.14296441537405
$ws.Range("L7").Value2 = 0.0406192295686798
$ws.Range("B8").Value2 = 0.123520085934971
$ws.Range("C8").Value2 = 0.116423087925096
$ws.Range("D8").Value2 = 0.145048290059828
$ws.Range("E8").Value2 = 0.0259849800701356
$ws.Range("F8").Value2 = 0.150124247211813
$ws.Range("G8").Value2 = 0.0411385001482822
$ws.Range("H8").Value2 = 0.120459163349451
$ws.Range("I8").Value2 = 0.213078974573974
$ws.Range("J8").Value2 = 0.0408066655072365
$ws.Range("K8").Value2 = 0.0171711263334772
$ws.Range("L8").Value2 = 0.0425388041985455
$ws.Range("M8").Value2 = 0.188101401357083
$ws.Range("N8").Value2 = 0.201846451513117
$ws.Range("B9").Value2 = 0.123173868549124
$ws.Range("C9").Value2 = 0.10975621589324
$ws.Range("D9").Value2 = 0.0752824629903547
$ws.Range("E9").Value2 = 0.0788480185772115
$ws.Range("F9").Value2 = 0.12839111657431
$ws.Range("G9").Value2 = 0.0551413945323128
$ws.Range("H9").Value2 = 0.133706009297009
$ws.Range("I9").Value2 = 0.157768426282657
$ws.Range("J9").Value2 = 0.0551564177571546
$ws.Range("K9").Value2 = 0.0764199748656044
$ws.Range("L9").Value2 = 0.0490690003768072
$ws.Range("M9").Value2 = 0.113148347477788
$ws.Range("N9").Value2 = 0.185150157064038
$ws.Range("B10").Value2 = 0.108516161624555
$ws.Range("L10").Value2 = 0.0547249303757871
$ws.Range("B11").Value2 = 0.108054915619012
$ws.Range("C11").Value2 = 0.0671024241823383
$ws.Range("D11").Value2 = 0.0907094935471546
$ws.Range("E11").Value2 = 0.0563285510537544
$ws.Range("F11").Value2 = 0.0432483434768305
$ws.Range("G11").Value2 = 0.0308288925637526
$ws.Range("H11").Value2 = 0.081452689714342
$ws.Range("I11").Value2 = 0.0873940812101499
$ws.Range("J11").Value2 = 0.0631680841066284
$ws.Range("K11").Value2 = 0.155618854206164
$ws.Range("L11").Value2 = 0.26382606513072
$ws.Range("M11").Value2 = 0.0709578059931298
$ws.Range("N11").Value2 = 0.0719106694806072
$ws.Range("B12").Value2 = 0.105897273746632
$ws.Range("C12").Value2 = 0.0977448568617962
$ws.Range("D12").Value2 = 0.107748561410452
$ws.Range("E12").Value2 = 0.118838926880645
$ws.Range("F12").Value2 = 0.0871212802155489
$ws.Range("G12").Value2 = 0.104562836800394
$ws.Range("H12").Value2 = 0.0803608597937457
$ws.Range("I12").Value2 = 0.062809296463231
$ws.Range("J12").Value2 = 0.193202005944839
$ws.Range("K12").Value2 = 0.12061045749135
$ws.Range("L12").Value2 = 0.271808786699362
$ws.Range("M12").Value2 = 0.035722429855013
$ws.Range("N12").Value2 = 0.0531058407190556
$ws.Range("B13").Value2 = 0.074783280562972
$ws.Range("L13").Value2 = 0.0114558411970786
$ws.Range("B14").Value2 = 0.0573135385366308
$ws.Range("L14").Value2 = 0.0382615295733116
$ws.Range("B15").Value2 = 0.0548094909108072
$ws.Range("L15").Value2 = 0.115194950755808
$ws.Range("B16").Value2 = 0.0524527620637507
$ws.Range("L16").Value2 = 0.0848437475253124
$ws.Range("B17").Value2 = 0.0483232859942362
$ws.Range("C17").Value2 = 0.020107200342392
$ws.Range("D17").Value2 = 0.0191072298135684
$ws.Range("E17").Value2 = 0.0147445422463436
$ws.Range("F17").Value2 = 0.0422960057733754
$ws.Range("G17").Value2 = 0.0138812498374743
$ws.Range("H17").Value2 = 0
$ws.Range("I17").Value2 = 0.0238569653014895
$ws.Range("J17").Value2 = 0.016704071214546
$ws.Range("K17").Value2 = 0.0432944191285272
$ws.Range("L17").Value2 = 0.114376058589707
$ws.Range("M17").Value2 = 0.177940248053361
$ws.Range("N17").Value2 = 0.0447466308579002
$ws.Range("B18").Value2 = 0.0450924273488087
$ws.Range("C18").Value2 = 0.0318956347139389
$ws.Range("D18").Value2 = 0.0347569102343155
$ws.Range("E18").Value2 = 0.0227949352075646
$ws.Range("F18").Value2 = 0.0503326400232044
$ws.Range("G18").Value2 = 0.0398434107874001
$ws.Range("H18").Value2 = 0.0269444936151781
$ws.Range("I18").Value2 = 0.0130829252163652
$ws.Range("J18").Value2 = 0.0746011086790753
$ws.Range("K18").Value2 = 0.0587829032453478
$ws.Range("L18").Value2 = 0.0423199026190579
$ws.Range("M18").Value2 = 0.063715983500485
$ws.Range("N18").Value2 = 0.0518581935554064
$ws.Range("B19").Value2 = 0.0419122271591407
$ws.Range("C19").Value2 = 0.0159273771810526
$ws.Range("D19").Value2 = 0.00431485925509307
$ws.Range("E19").Value2 = 0.016641649524601
$ws.Range("F19").Value2 = 0.0177719932615889
$ws.Range("G19").Value2 = 0.0223840174246049
$ws.Range("H19").Value2 = 0.00528653334412038
$ws.Range("I19").Value2 = 0.0263976937533577
$ws.Range("J19").Value2 = 0.0110943078154441
$ws.Range("K19").Value2 = 0.00363309988422034
$ws.Range("L19").Value2 = 0
$ws.Range("M19").Value2 = 0
$ws.Range("N19").Value2 = 0.10545839437447
$ws.Range("B20").Value2 = 0.0403629758493454
$ws.Range("C20").Value2 = 0.0465839083523768
$ws.Range("D20").Value2 = 0.0303596134606511
$ws.Range("E20").Value2 = 0.067759608960212
$ws.Range("F20").Value2 = 0.0174913878758809
$ws.Range("G20").Value2 = 0.036005597790544
$ws.Range("H20").Value2 = 0.0301003740322274
$ws.Range("I20").Value2 = 0.0822581165295588
$ws.Range("J20").Value2 = 0.035594242190962
$ws.Range("K20").Value2 = 0.0260377303001874
$ws.Range("L20").Value2 = 0.00629614863032136
$ws.Range("M20").Value2 = 0.0116800076603586
$ws.Range("N20").Value2 = 0.0541361598107681
$ws.Range("B21").Value2 = 0.0326044250127659
$ws.Range("C21").Value2 = 0.0337159523942901
$ws.Range("D21").Value2 = 0.0350914726918312
$ws.Range("E21").Value2 = 0.019411861832518
$ws.Range("F21").Value2 = 0.017089497696424
$ws.Range("G21").Value2 = 0.00930785373822185
$ws.Range("H21").Value2 = 0.0669189269466735
$ws.Range("I21").Value2 = 0.0458679954908412
$ws.Range("J21").Value2 = 0.0169749940170554
$ws.Range("K21").Value2 = 0.00761920287517245
$ws.Range("L21").Value2 = 0.0573762896815998
$ws.Range("M21").Value2 = 0
$ws.Range("N21").Value2 = 0.0325956882961598
$ws.Range("B22").Value2 = 0.0314647462468635
$ws.Range("L22").Value2 = 0.0406227341544564
$ws.Range("B23").Value2 = 0.0312921717367679
$ws.Range("C23").Value2 = 0.017031276019717
$ws.Range("D23").Value2 = 0.0142180700661177
$ws.Range("E23").Value2 = 0.0350788999893401
$ws.Range("F23").Value2 = 0.00468058621332618
$ws.Range("G23").Value2 = 0
$ws.Range("H23").Value2 = 0.00830775886317187
$ws.Range("I23").Value2 = 0.0275557774206001
$ws.Range("J23").Value2 = 0
$ws.Range("K23").Value2 = 0.0082594908632089
$ws.Range("L23").Value2 = 0.0639944335058778
$ws.Range("M23").Value2 = 0.0372946625769812
$ws.Range("N23").Value2 = 0.0436983227225689
$ws.Range("B24").Value2 = 0.0277108854658857
$ws.Range("C24").Value2 = 0.0246002728483537
$ws.Range("D24").Value2 = 0.0188877599120833
$ws.Range("E24").Value2 = 0.0619870811812808
$ws.Range("F24").Value2 = 0.00819962221981932
$ws.Range("G24").Value2 = 0.00578913863703181
$ws.Range("H24").Value2 = 0.0169166175940457
$ws.Range("I24").Value2 = 0.0161053819456598
$ws.Range("J24").Value2 = 0.0191112586923636
$ws.Range("K24").Value2 = 0.0558648964707215
$ws.Range("L24").Value2 = 0.0273815458655608
$ws.Range("M24").Value2 = 0.0346828509910579
$ws.Range("N24").Value2 = 0.0198152593563636
$ws.Range("B25").Value2 = 0.0197176026547576
$ws.Range("L25").Value2 = 0.0366245499224145
$ws.Range("B26").Value2 = 0.0153318121650389
$ws.Range("L26").Value2 = 0
$ws.Range("B27").Value2 = 0.0150792992577955
$ws.Range("L27").Value2 = 0.0123042531636731
$ws.Range("B28").Value2 = 0.00936128373356803
$ws.Range("L28").Value2 = 0.00591440312288486
